# Update gh-pages to output generated at 456a3b4
# Applies updated '想去人数' (want-to-go count) figures across all sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5906
$ws.Range("F5").Value = 5906
$ws.Range("F7").Value = 2956
$ws.Range("F8").Value = 1272
$ws.Range("F12").Value = 698
$ws.Range("F13").Value = 212
$ws.Range("F14").Value = 4272
$ws.Range("F15").Value = 4272
$ws.Range("F18").Value = 105
$ws.Range("F22").Value = 6464
$ws.Range("F23").Value = 6465
$ws.Range("F24").Value = 226
$ws.Range("F25").Value = 95
$ws.Range("F27").Value = 447
$ws.Range("F28").Value = 1224
$ws.Range("F30").Value = 6242
$ws.Range("F31").Value = 1623
$ws.Range("F33").Value = 1857
$ws.Range("F34").Value = 5947
$ws.Range("F39").Value = 391
$ws.Range("F40").Value = 4069
$ws.Range("F41").Value = 5
$ws.Range("F42").Value = 186
$ws.Range("F50").Value = 310
$ws.Range("F51").Value = 2038

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 192
$ws.Range("F4").Value = 28
$ws.Range("F5").Value = 98
$ws.Range("F6").Value = 27
$ws.Range("F11").Value = 16

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1410

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1410
$ws.Range("F4").Value = 5907
$ws.Range("F5").Value = 5907
$ws.Range("F7").Value = 2956
$ws.Range("F8").Value = 1272
$ws.Range("F12").Value = 192
$ws.Range("F13").Value = 212
$ws.Range("F14").Value = 4272
$ws.Range("F15").Value = 4272
$ws.Range("F18").Value = 105
$ws.Range("F22").Value = 6466
$ws.Range("F23").Value = 6466
$ws.Range("F24").Value = 226
$ws.Range("F25").Value = 95
$ws.Range("F26").Value = 447
$ws.Range("F27").Value = 1224
$ws.Range("F28").Value = 98
$ws.Range("F29").Value = 6242
$ws.Range("F30").Value = 1623
$ws.Range("F31").Value = 27
$ws.Range("F33").Value = 1857
$ws.Range("F35").Value = 5947
$ws.Range("F40").Value = 391
$ws.Range("F41").Value = 4069
$ws.Range("F42").Value = 186
$ws.Range("F51").Value = 310
$ws.Range("F52").Value = 16

